# VandV Comparison.xlsx -- removed unused variable NEUTRAL from CFAST data common block
# Updates several J/K column values on Sheet1 (columns computed from the new
# common-block layout), replacing some static K-column offsets with
# formulas derived from the new (sigma-M)/2 relationship, and clears two
# now-unused K-column zeros. Also updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("J3").Value = 1.1100000000000001
$ws.Range("K3").Formula = "=0.44/2"

# Row 4
$ws.Range("J4").Value = 1.01
$ws.Range("K4").Formula = "=0.32/2"

# Row 5
$ws.Range("J5").Value = 1.25
$ws.Range("K5").Formula = "=0.53/2"

# Row 6
$ws.Range("K6").Formula = "=0.42/2"

# Row 8
$ws.Range("J8").Value = 1.03
$ws.Range("K8").Formula = "=0.63/2"

# Row 9
$ws.Range("K9").Formula = "=0.56/2"

# Row 10 / 11 - K column no longer holds a literal 0; clear the cell
$ws.Range("K10").ClearContents()
$ws.Range("K11").ClearContents()

# Row 13
$ws.Range("K13").Formula = "=1.29/2"

# Row 15
$ws.Range("J15").Value = 0.99
$ws.Range("K15").Formula = "=0.99/2"

# Update the active selection on the sheet (was F12, now J1:L1)
$ws.Range("J1:L1").Select()
